# Lecture 3 slides: update the "formal languages" bullet on the
# "Formal Languages" slide (sldId 258 -> 4th slide in the deck) to a
# more concise wording.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Paragraphs(2).Text = "Formal languages follow strict rule"
